$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update Fase 3: Garantia de Calidad y Depuracion entries
$ws.Range("C28").Value = 8
$ws.Range("D28").Value = 5
$ws.Range("D29").Value = 2
$ws.Range("D30").Value = 4

# Update Fase 4: Documentacion del proyecto entries
$ws.Range("D35").Value = 4

# Reflect the final selection the author left the sheet on
$ws.Range("E38").Select()
